$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 179
$ws.Cells.Item(179, 2).Value = 6992695
$ws.Cells.Item(179, 5).Value = "Muang Thong United"
$ws.Cells.Item(179, 6).Value = "Uthai Thani FC"
$ws.Cells.Item(179, 7).Value = 5
$ws.Cells.Item(179, 8).Value = 2
$ws.Cells.Item(179, 9).Value = 0
$ws.Cells.Item(179, 10).Value = 2
$ws.Cells.Item(179, 11).Value = "H"
$ws.Cells.Item(179, 12).Value = 2.1
$ws.Cells.Item(179, 13).Value = 3.75
$ws.Cells.Item(179, 14).Value = 2.7
$ws.Cells.Item(179, 15).Value = 1.95
$ws.Cells.Item(179, 16).Value = 3.8
$ws.Cells.Item(179, 17).Value = 2.9
$ws.Cells.Item(179, 18).Value = -0.25
$ws.Cells.Item(179, 19).Value = 1.8
$ws.Cells.Item(179, 20).Value = 2
$ws.Cells.Item(179, 21).Value = 3
$ws.Cells.Item(179, 22).Value = 1.825
$ws.Cells.Item(179, 23).Value = 1.975
$ws.Cells.Item(179, 24).Value = 0.95
$ws.Cells.Item(179, 25).Value = -1
$ws.Cells.Item(179, 26).Value = -1
$ws.Cells.Item(179, 27).Value = 0.8
$ws.Cells.Item(179, 28).Value = -1
$ws.Cells.Item(179, 29).Value = 0.825
$ws.Cells.Item(179, 30).Value = -1

# Row 180
$ws.Cells.Item(180, 2).Value = 8026714
$ws.Cells.Item(180, 5).Value = "BG Pathum United"
$ws.Cells.Item(180, 6).Value = "Buriram United"
$ws.Cells.Item(180, 7).Value = 1
$ws.Cells.Item(180, 8).Value = 1
$ws.Cells.Item(180, 9).Value = 0
$ws.Cells.Item(180, 10).Value = 1
$ws.Cells.Item(180, 11).Value = "D"
$ws.Cells.Item(180, 12).Value = 3
$ws.Cells.Item(180, 13).Value = 3.6
$ws.Cells.Item(180, 14).Value = 2
$ws.Cells.Item(180, 15).Value = 3.1
$ws.Cells.Item(180, 16).Value = 3.75
$ws.Cells.Item(180, 17).Value = 1.95
$ws.Cells.Item(180, 18).Value = 0.5
$ws.Cells.Item(180, 19).Value = 1.825
$ws.Cells.Item(180, 20).Value = 1.975
$ws.Cells.Item(180, 21).Value = 2.75
$ws.Cells.Item(180, 22).Value = 1.85
$ws.Cells.Item(180, 23).Value = 1.95
$ws.Cells.Item(180, 24).Value = -1
$ws.Cells.Item(180, 25).Value = 2.75
$ws.Cells.Item(180, 26).Value = -1
$ws.Cells.Item(180, 27).Value = 0.825
$ws.Cells.Item(180, 28).Value = -1
$ws.Cells.Item(180, 29).Value = -1
$ws.Cells.Item(180, 30).Value = 0.95

# Row 225
$ws.Cells.Item(225, 2).Value = 6992741
$ws.Cells.Item(225, 5).Value = "Muang Thong United"
$ws.Cells.Item(225, 6).Value = "Chonburi"
$ws.Cells.Item(225, 7).Value = 6
$ws.Cells.Item(225, 8).Value = 0
$ws.Cells.Item(225, 9).Value = 1
$ws.Cells.Item(225, 10).Value = 0
$ws.Cells.Item(225, 11).Value = "H"
$ws.Cells.Item(225, 12).Value = 1.75
$ws.Cells.Item(225, 13).Value = 4
$ws.Cells.Item(225, 14).Value = 3.4
$ws.Cells.Item(225, 15).Value = 1.95
$ws.Cells.Item(225, 16).Value = 4
$ws.Cells.Item(225, 17).Value = 2.9
$ws.Cells.Item(225, 18).Value = -0.25
$ws.Cells.Item(225, 19).Value = 1.75
$ws.Cells.Item(225, 20).Value = 1.95
$ws.Cells.Item(225, 21).Value = 3.25
$ws.Cells.Item(225, 22).Value = 1.875
$ws.Cells.Item(225, 23).Value = 1.925
$ws.Cells.Item(225, 24).Value = 0.95
$ws.Cells.Item(225, 25).Value = -1
$ws.Cells.Item(225, 26).Value = -1
$ws.Cells.Item(225, 27).Value = 0.75
$ws.Cells.Item(225, 28).Value = -1
$ws.Cells.Item(225, 29).Value = 0.875
$ws.Cells.Item(225, 30).Value = -1

# Row 226
$ws.Cells.Item(226, 2).Value = 6992738
$ws.Cells.Item(226, 5).Value = "Prachuap FC"
$ws.Cells.Item(226, 6).Value = "Chiangrai Utd"
$ws.Cells.Item(226, 7).Value = 1
$ws.Cells.Item(226, 8).Value = 0
$ws.Cells.Item(226, 9).Value = 1
$ws.Cells.Item(226, 10).Value = 0
$ws.Cells.Item(226, 11).Value = "H"
$ws.Cells.Item(226, 12).Value = 1.666
$ws.Cells.Item(226, 13).Value = 3.5
$ws.Cells.Item(226, 14).Value = 4.5
$ws.Cells.Item(226, 15).Value = 1.48
$ws.Cells.Item(226, 16).Value = 4
$ws.Cells.Item(226, 17).Value = 5.75
$ws.Cells.Item(226, 18).Value = -1.25
$ws.Cells.Item(226, 19).Value = 2.025
$ws.Cells.Item(226, 20).Value = 1.775
$ws.Cells.Item(226, 21).Value = 2.75
$ws.Cells.Item(226, 22).Value = 1.975
$ws.Cells.Item(226, 23).Value = 1.825
$ws.Cells.Item(226, 24).Value = 0.48
$ws.Cells.Item(226, 25).Value = -1
$ws.Cells.Item(226, 26).Value = -1
$ws.Cells.Item(226, 27).Value = -0.5
$ws.Cells.Item(226, 28).Value = 0.3875
$ws.Cells.Item(226, 29).Value = -1
$ws.Cells.Item(226, 30).Value = 0.825

# Row 232
$ws.Cells.Item(232, 2).Value = 6992748
$ws.Cells.Item(232, 5).Value = "Port FC"
$ws.Cells.Item(232, 6).Value = "Nakhon Pathom FC"
$ws.Cells.Item(232, 7).Value = 6
$ws.Cells.Item(232, 8).Value = 0
$ws.Cells.Item(232, 11).Value = "H"
$ws.Cells.Item(232, 12).Value = 1.3
$ws.Cells.Item(232, 13).Value = 5.25
$ws.Cells.Item(232, 14).Value = 6.25
$ws.Cells.Item(232, 15).Value = 1.363
$ws.Cells.Item(232, 16).Value = 5
$ws.Cells.Item(232, 17).Value = 5.75
$ws.Cells.Item(232, 18).Value = -1.5
$ws.Cells.Item(232, 19).Value = 1.875
$ws.Cells.Item(232, 20).Value = 1.925
$ws.Cells.Item(232, 21).Value = 3.5
$ws.Cells.Item(232, 22).Value = 1.925
$ws.Cells.Item(232, 23).Value = 1.775
$ws.Cells.Item(232, 24).Value = 0.363
$ws.Cells.Item(232, 25).Value = -1
$ws.Cells.Item(232, 26).Value = -1
$ws.Cells.Item(232, 27).Value = 0.875
$ws.Cells.Item(232, 28).Value = -1
$ws.Cells.Item(232, 29).Value = 0.925
$ws.Cells.Item(232, 30).Value = -1

# Row 233
$ws.Cells.Item(233, 2).Value = 6992746
$ws.Cells.Item(233, 5).Value = "Chiangrai Utd"
$ws.Cells.Item(233, 6).Value = "Muang Thong United"
$ws.Cells.Item(233, 7).Value = 2
$ws.Cells.Item(233, 8).Value = 3
$ws.Cells.Item(233, 11).Value = "A"
$ws.Cells.Item(233, 12).Value = 3
$ws.Cells.Item(233, 13).Value = 3.4
$ws.Cells.Item(233, 14).Value = 2.1
$ws.Cells.Item(233, 15).Value = 2.5
$ws.Cells.Item(233, 16).Value = 3.4
$ws.Cells.Item(233, 17).Value = 2.4
$ws.Cells.Item(233, 18).Value = 0
$ws.Cells.Item(233, 19).Value = 1.975
$ws.Cells.Item(233, 20).Value = 1.825
$ws.Cells.Item(233, 21).Value = 2.75
$ws.Cells.Item(233, 22).Value = 1.825
$ws.Cells.Item(233, 23).Value = 1.975
$ws.Cells.Item(233, 24).Value = -1
$ws.Cells.Item(233, 25).Value = -1
$ws.Cells.Item(233, 26).Value = 1.4
$ws.Cells.Item(233, 27).Value = -1
$ws.Cells.Item(233, 28).Value = 0.825
$ws.Cells.Item(233, 29).Value = 0.825
$ws.Cells.Item(233, 30).Value = -1

# Row 234
$ws.Cells.Item(234, 2).Value = 6992745
$ws.Cells.Item(234, 5).Value = "Sukhothai FC"
$ws.Cells.Item(234, 6).Value = "Lamphun Warrior FC"
$ws.Cells.Item(234, 7).Value = 0
$ws.Cells.Item(234, 8).Value = 3
$ws.Cells.Item(234, 11).Value = "A"
$ws.Cells.Item(234, 12).Value = 2.3
$ws.Cells.Item(234, 13).Value = 3.25
$ws.Cells.Item(234, 14).Value = 2.875
$ws.Cells.Item(234, 15).Value = 3
$ws.Cells.Item(234, 16).Value = 3.3
$ws.Cells.Item(234, 17).Value = 2.15
$ws.Cells.Item(234, 18).Value = 0.25
$ws.Cells.Item(234, 19).Value = 1.9
$ws.Cells.Item(234, 20).Value = 1.9
$ws.Cells.Item(234, 21).Value = 2.75
$ws.Cells.Item(234, 22).Value = 1.85
$ws.Cells.Item(234, 23).Value = 1.95
$ws.Cells.Item(234, 24).Value = -1
$ws.Cells.Item(234, 25).Value = -1
$ws.Cells.Item(234, 26).Value = 1.15
$ws.Cells.Item(234, 27).Value = -1
$ws.Cells.Item(234, 28).Value = 0.8999999999999999
$ws.Cells.Item(234, 29).Value = 0.425
$ws.Cells.Item(234, 30).Value = -0.5

# Row 235
$ws.Cells.Item(235, 2).Value = 6992750
$ws.Cells.Item(235, 5).Value = "Bangkok United"
$ws.Cells.Item(235, 6).Value = "Uthai Thani FC"
$ws.Cells.Item(235, 7).Value = 3
$ws.Cells.Item(235, 8).Value = 0
$ws.Cells.Item(235, 11).Value = "H"
$ws.Cells.Item(235, 12).Value = 1.3
$ws.Cells.Item(235, 13).Value = 5.5
$ws.Cells.Item(235, 14).Value = 7
$ws.Cells.Item(235, 15).Value = 1.285
$ws.Cells.Item(235, 16).Value = 6
$ws.Cells.Item(235, 17).Value = 7
$ws.Cells.Item(235, 18).Value = -1.75
$ws.Cells.Item(235, 19).Value = 1.875
$ws.Cells.Item(235, 20).Value = 1.925
$ws.Cells.Item(235, 21).Value = 3.5
$ws.Cells.Item(235, 22).Value = 1.975
$ws.Cells.Item(235, 23).Value = 1.825
$ws.Cells.Item(235, 24).Value = 0.2849999999999999
$ws.Cells.Item(235, 25).Value = -1
$ws.Cells.Item(235, 26).Value = -1
$ws.Cells.Item(235, 27).Value = 0.875
$ws.Cells.Item(235, 28).Value = -1
$ws.Cells.Item(235, 29).Value = -1
$ws.Cells.Item(235, 30).Value = 0.825

# Row 237
$ws.Cells.Item(237, 2).Value = 6992342
$ws.Cells.Item(237, 5).Value = "BG Pathum United"
$ws.Cells.Item(237, 6).Value = "Police Tero FC"
$ws.Cells.Item(237, 7).Value = 4
$ws.Cells.Item(237, 8).Value = 2
$ws.Cells.Item(237, 11).Value = "H"
$ws.Cells.Item(237, 12).Value = 1.25
$ws.Cells.Item(237, 13).Value = 6
$ws.Cells.Item(237, 14).Value = 7.5
$ws.Cells.Item(237, 15).Value = 1.25
$ws.Cells.Item(237, 16).Value = 6
$ws.Cells.Item(237, 17).Value = 7.5
$ws.Cells.Item(237, 18).Value = -2
$ws.Cells.Item(237, 19).Value = 1.95
$ws.Cells.Item(237, 20).Value = 1.85
$ws.Cells.Item(237, 21).Value = 3.75
$ws.Cells.Item(237, 22).Value = 1.85
$ws.Cells.Item(237, 23).Value = 1.95
$ws.Cells.Item(237, 24).Value = 0.25
$ws.Cells.Item(237, 25).Value = -1
$ws.Cells.Item(237, 26).Value = -1
$ws.Cells.Item(237, 27).Value = 0
$ws.Cells.Item(237, 28).Value = 0
$ws.Cells.Item(237, 29).Value = 0.8500000000000001
$ws.Cells.Item(237, 30).Value = -1

# Row 238
$ws.Cells.Item(238, 2).Value = 6992744
$ws.Cells.Item(238, 5).Value = "Chonburi"
$ws.Cells.Item(238, 6).Value = "Trat FC"
$ws.Cells.Item(238, 7).Value = 3
$ws.Cells.Item(238, 8).Value = 2
$ws.Cells.Item(238, 11).Value = "H"
$ws.Cells.Item(238, 12).Value = 1.615
$ws.Cells.Item(238, 13).Value = 3.9
$ws.Cells.Item(238, 14).Value = 4.5
$ws.Cells.Item(238, 15).Value = 1.571
$ws.Cells.Item(238, 16).Value = 4.1
$ws.Cells.Item(238, 17).Value = 4.333
$ws.Cells.Item(238, 18).Value = -1
$ws.Cells.Item(238, 19).Value = 1.95
$ws.Cells.Item(238, 20).Value = 1.85
$ws.Cells.Item(238, 21).Value = 3.5
$ws.Cells.Item(238, 22).Value = 1.975
$ws.Cells.Item(238, 23).Value = 1.825
$ws.Cells.Item(238, 24).Value = 0.571
$ws.Cells.Item(238, 25).Value = -1
$ws.Cells.Item(238, 26).Value = -1
$ws.Cells.Item(238, 27).Value = 0
$ws.Cells.Item(238, 28).Value = 0
$ws.Cells.Item(238, 29).Value = 0.9750000000000001
$ws.Cells.Item(238, 30).Value = -1
